$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.908.26'
$ws.Range("E2").Value = '  -0.33%  '

$ws.Range("D3").Value = '3.846.13'
$ws.Range("E3").Value = '  +0.92%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '699.44'
$ws.Range("E5").Value = '  -1.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.91'
$ws.Range("E6").Value = '  -0.92%  '

$ws.Range("D7").Value = '3.844.30'
$ws.Range("E7").Value = '  +0.91%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("E9").Value = '  -0.57%  '

$ws.Range("E10").Value = '  -1.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.28'
$ws.Range("E11").Value = '  -1.96%  '

$ws.Range("E12").Value = '  -0.90%  '

$ws.Range("E13").Value = '  -0.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.23'
$ws.Range("E14").Value = '  -0.40%  '

$ws.Range("E15").Value = '  +0.98%  '

$ws.Range("D16").Value = '3.842.93'
$ws.Range("E16").Value = '  +0.75%  '

$ws.Range("D17").Value = '70.911.51'
$ws.Range("E17").Value = '  -0.35%  '

$ws.Range("E18").Value = '  -1.32%  '

$ws.Range("E19").Value = '  +0.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.43'
$ws.Range("E20").Value = '  -3.28%  '

$ws.Range("E21").Value = '  -4.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '493.68'
$ws.Range("E22").Value = '  +1.95%  '

$ws.Range("E23").Value = '  -0.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.55'
$ws.Range("E24").Value = '  +0.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000148'
$ws.Range("E25").Value = '  +0.80%  '

$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.56'
$ws.Range("E26").Value = '  -1.27%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.15'
$ws.Range("E27").Value = '  -3.14%  '

$ws.Range("E28").Value = '  -4.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.12'
$ws.Range("E30").Value = '  +1.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.51'
$ws.Range("E31").Value = '  -1.40%  '

$ws.Range("E32").Value = '  -2.20%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.45'
$ws.Range("E33").Value = '  -1.00%  '

$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.182'
$ws.Range("E34").Value = '  +1.43%  '

$ws.Range("D35").Value = '3.801.22'
$ws.Range("E35").Value = '  +1.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.14'
$ws.Range("E36").Value = '  -1.61%  '

$ws.Range("E37").Value = '  +0.03%  '

$ws.Range("E38").Value = '  -1.10%  '

$ws.Range("E39").Value = '  +5.14%  '

$ws.Range("E40").Value = '  +6.79%  '

$ws.Range("E41").Value = '  -0.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.34'
$ws.Range("E42").Value = '  -5.05%  '

$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000313'
$ws.Range("E45").Value = '  -8.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '163.55'
$ws.Range("E46").Value = '  +1.89%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '48.70'
$ws.Range("E47").Value = '  -1.57%  '

$ws.Range("E48").Value = '  -1.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.63'
$ws.Range("E49").Value = '  +0.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.30'
$ws.Range("E50").Value = '  -5.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.36'
$ws.Range("E51").Value = '  -4.97%  '
